$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2048.182
$ws.Range("I28").Value = 1392.2222
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 1392.2222
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = -907.2221999999999
$ws.Range("N28").Value = -5970
$ws.Range("H62").Value = 3503.0588
$ws.Range("I62").Value = 3466.5386
$ws.Range("J62").Value = 3621.75
$ws.Range("K62").Value = 3466.5386
$ws.Range("L62").Value = 3621.75
$ws.Range("M62").Value = -2842.5386
$ws.Range("N62").Value = -4869.75
$ws.Range("H65").Value = 3503.0588
$ws.Range("I65").Value = 3466.5386
$ws.Range("J65").Value = 3621.75
$ws.Range("K65").Value = 17332.693
$ws.Range("L65").Value = 18108.75
$ws.Range("M65").Value = -14212.693
$ws.Range("N65").Value = -24348.75
$ws.Range("H98").Value = 925
$ws.Range("I98").Value = 945.2381
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 945.2381
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 552.7619
$ws.Range("N98").Value = -3496
$ws.Range("H106").Value = 4183.7646
$ws.Range("I106").Value = 3152.0715
$ws.Range("K106").Value = 3152.0715
$ws.Range("M106").Value = -2521.0715
$ws.Range("H107").Value = 511
$ws.Range("I107").Value = 491.25
$ws.Range("J107").Value = 590
$ws.Range("K107").Value = 491.25
$ws.Range("L107").Value = 590
$ws.Range("M107").Value = 1428.75
$ws.Range("N107").Value = -4430
$ws.Range("H122").Value = 925
$ws.Range("I122").Value = 945.2381
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 2835.7143
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -385.7143000000001
$ws.Range("N122").Value = -6400
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H126").Value = 35000
$ws.Range("J126").Value = 35000
$ws.Range("L126").Value = 35000
$ws.Range("N126").Value = -44880
$ws.Range("H137").Value = 784209.75
$ws.Range("I137").Value = 1854.9375
$ws.Range("J137").Value = 1426141.9
$ws.Range("K137").Value = 5564.8125
$ws.Range("L137").Value = 4278425.699999999
$ws.Range("M137").Value = -3014.8125
$ws.Range("N137").Value = -4283525.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2421.2
$ws.Range("I45").Value = 2494.5715
$ws.Range("K45").Value = 2494.5715
$ws.Range("M45").Value = -2117.5715
$ws.Range("H110").Value = 2223.3
$ws.Range("I110").Value = 2223.3
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2223.3
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -178.3000000000002
$ws.Range("N110").Value = ""
$ws.Range("H135").Value = 27362.375
$ws.Range("J135").Value = 27362.375
$ws.Range("L135").Value = 27362.375
$ws.Range("N135").Value = -37502.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 50250
$ws.Range("J39").Value = 50250
$ws.Range("L39").Value = 50250
$ws.Range("N39").Value = -51028

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""
$ws.Range("H25").Value = 3011
$ws.Range("I25").Value = 3011
$ws.Range("K25").Value = 3011
$ws.Range("M25").Value = -2837
$ws.Range("H31").Value = 7813652
$ws.Range("I31").Value = 743.6957
$ws.Range("J31").Value = 12196503
$ws.Range("K31").Value = 743.6957
$ws.Range("L31").Value = 12196503
$ws.Range("M31").Value = -448.6957
$ws.Range("N31").Value = -12197093
$ws.Range("H34").Value = 7813652
$ws.Range("I34").Value = 743.6957
$ws.Range("J34").Value = 12196503
$ws.Range("K34").Value = 743.6957
$ws.Range("L34").Value = 12196503
$ws.Range("M34").Value = -541.6957
$ws.Range("N34").Value = -12196907
$ws.Range("H86").Value = 4726.067
$ws.Range("I86").Value = 2294.75
$ws.Range("K86").Value = 2294.75
$ws.Range("M86").Value = -1171.75
$ws.Range("H89").Value = 4726.067
$ws.Range("I89").Value = 2294.75
$ws.Range("K89").Value = 11473.75
$ws.Range("M89").Value = -5857.75
$ws.Range("H99").Value = 1754.8966
$ws.Range("I99").Value = 1823.238
$ws.Range("J99").Value = 1575.5
$ws.Range("K99").Value = 1823.238
$ws.Range("L99").Value = 1575.5
$ws.Range("M99").Value = -325.2380000000001
$ws.Range("N99").Value = -4571.5
$ws.Range("H126").Value = 1754.8966
$ws.Range("I126").Value = 1823.238
$ws.Range("J126").Value = 1575.5
$ws.Range("K126").Value = 5469.714
$ws.Range("L126").Value = 4726.5
$ws.Range("M126").Value = -2999.714
$ws.Range("N126").Value = -9666.5
$ws.Range("H132").Value = 3811.32
$ws.Range("I132").Value = 3768.0588
$ws.Range("J132").Value = 3903.25
$ws.Range("K132").Value = 11304.1764
$ws.Range("L132").Value = 11709.75
$ws.Range("M132").Value = -8774.1764
$ws.Range("N132").Value = -16769.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 19299044
$ws.Range("I5").Value = 20000568
$ws.Range("J5").Value = 18519572
$ws.Range("K5").Value = 60001704
$ws.Range("L5").Value = 55558716
$ws.Range("M5").Value = -60001592
$ws.Range("N5").Value = -55558940
$ws.Range("H80").Value = 1460.8
$ws.Range("I80").Value = 1802
$ws.Range("J80").Value = 1233.3334
$ws.Range("K80").Value = 5406
$ws.Range("L80").Value = 3700.0002
$ws.Range("M80").Value = -4470
$ws.Range("N80").Value = -5572.0002
$ws.Range("H83").Value = 1460.8
$ws.Range("I83").Value = 1802
$ws.Range("J83").Value = 1233.3334
$ws.Range("K83").Value = 16218
$ws.Range("L83").Value = 11100.0006
$ws.Range("M83").Value = -11538
$ws.Range("N83").Value = -20460.0006
$ws.Range("H131").Value = 880.6799999999999
$ws.Range("J131").Value = 895.7083
$ws.Range("L131").Value = 2687.1249
$ws.Range("N131").Value = -12767.1249
$ws.Range("H135").Value = 19299044
$ws.Range("I135").Value = 20000568
$ws.Range("J135").Value = 18519572
$ws.Range("K135").Value = 180005112
$ws.Range("L135").Value = 166676148
$ws.Range("M135").Value = -180002577
$ws.Range("N135").Value = -166681218

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2896.8333
$ws.Range("I31").Value = 1476.2
$ws.Range("K31").Value = 1476.2
$ws.Range("M31").Value = -1184.2
$ws.Range("H37").Value = 2896.8333
$ws.Range("I37").Value = 1476.2
$ws.Range("K37").Value = 1476.2
$ws.Range("M37").Value = -1199.2
$ws.Range("H113").Value = 1481.3334
$ws.Range("I113").Value = 1391.5
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1391.5
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 778.5
$ws.Range("N113").Value = -6540
$ws.Range("H122").Value = 201018.6
$ws.Range("I122").Value = 201018.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 603055.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -600605.8
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2762.6667
$ws.Range("I7").Value = 2098.8572
$ws.Range("J7").Value = 3692
$ws.Range("K7").Value = 2098.8572
$ws.Range("L7").Value = 3692
$ws.Range("M7").Value = -1986.8572
$ws.Range("N7").Value = -3916
$ws.Range("H40").Value = 4290.8
$ws.Range("I40").Value = 4023.1
$ws.Range("J40").Value = 5897
$ws.Range("K40").Value = 4023.1
$ws.Range("L40").Value = 5897
$ws.Range("M40").Value = -3887.1
$ws.Range("N40").Value = -6169
$ws.Range("H122").Value = 2742
$ws.Range("I122").Value = 2617.2
$ws.Range("J122").Value = 3990
$ws.Range("K122").Value = 7851.599999999999
$ws.Range("L122").Value = 11970
$ws.Range("M122").Value = -5401.599999999999
$ws.Range("N122").Value = -16870
$ws.Range("H126").Value = 2762.6667
$ws.Range("I126").Value = 2098.8572
$ws.Range("J126").Value = 3692
$ws.Range("K126").Value = 6296.571599999999
$ws.Range("L126").Value = 11076
$ws.Range("M126").Value = -3826.571599999999
$ws.Range("N126").Value = -16016
$ws.Range("H132").Value = 4274.5137
$ws.Range("I132").Value = 3743.182
$ws.Range("J132").Value = 5053.8
$ws.Range("K132").Value = 11229.546
$ws.Range("L132").Value = 15161.4
$ws.Range("M132").Value = -8699.545999999998
$ws.Range("N132").Value = -20221.4
$ws.Range("H136").Value = 2090.5151
$ws.Range("I136").Value = 1938.9048
$ws.Range("J136").Value = 2355.8333
$ws.Range("K136").Value = 5816.7144
$ws.Range("L136").Value = 7067.499899999999
$ws.Range("M136").Value = -3266.7144
$ws.Range("N136").Value = -12167.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 7803
$ws.Range("I32").Value = 7008
$ws.Range("J32").Value = 8333
$ws.Range("K32").Value = 7008
$ws.Range("L32").Value = 8333
$ws.Range("M32").Value = -6691
$ws.Range("N32").Value = -8967
$ws.Range("H43").Value = 11475
$ws.Range("J43").Value = 11475
$ws.Range("L43").Value = 11475
$ws.Range("N43").Value = -11773
$ws.Range("H122").Value = 39089.03
$ws.Range("I122").Value = 1547.45
$ws.Range("J122").Value = 101658.336
$ws.Range("K122").Value = 4642.35
$ws.Range("L122").Value = 304975.008
$ws.Range("M122").Value = -2192.35
$ws.Range("N122").Value = -309875.008
$ws.Range("H126").Value = 83335250
$ws.Range("I126").Value = 1499
$ws.Range("J126").Value = 166669010
$ws.Range("K126").Value = 4497
$ws.Range("L126").Value = 500007030
$ws.Range("M126").Value = -2027
$ws.Range("N126").Value = -500011970
$ws.Range("H132").Value = 3657.9756
$ws.Range("I132").Value = 4400.4346
$ws.Range("J132").Value = 2709.2778
$ws.Range("K132").Value = 13201.3038
$ws.Range("L132").Value = 8127.8334
$ws.Range("M132").Value = -10671.3038
$ws.Range("N132").Value = -13187.8334

